$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 61.86
$ws.Range("E2").Value = -61.86

$ws.Range("D4").Value = 393.61
$ws.Range("E4").Value = 13329.73
$ws.Range("F4").Value = 0.02868179320777595
